$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph right before
# the section break. Expand it into five paragraphs: a bold+underlined
# "Gigs continued" sub-heading followed by four plain notes paragraphs.

# Index of the existing trailing empty paragraph - this won't shift as we
# append new paragraphs after it.
$firstIndex = $d.Paragraphs.Count

$last = $d.Paragraphs.Item($firstIndex).Range
$last.InsertParagraphAfter()
$d.Paragraphs.Item($firstIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstIndex + 3).Range.InsertParagraphAfter()

$texts = @(
    "Gigs continued",
    "Maybe use python using a form (I’ll create the python file but this depends on Sam), otherwise jquery",
    "A bunch of gigs that if possible that slide upwards after submission (jquery)",
    "Will need more images",
    "Use dates, name of artists and maybe venue"
)

for ($i = 0; $i -lt $texts.Length; $i++) {
    $p = $d.Paragraphs.Item($firstIndex + $i)
    $start = $p.Range.Start
    $insertPoint = $d.Range($start, $start)
    $insertPoint.InsertAfter($texts[$i])

    if ($i -eq 0) {
        $insertPoint.Font.Bold = $true
        $insertPoint.Font.BoldBi = $true
        $insertPoint.Font.Underline = 1
    }
}
